$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rows that got results filled in
$ws.Range("G84").Value = "Acierto"
$ws.Range("H84").Value = 0.53

$ws.Range("G85").Value = "Fallo"
$ws.Range("H85").Value = -1

$ws.Range("G98").Value = "Fallo"
$ws.Range("H98").Value = -1

$ws.Range("G99").Value = "Fallo"
$ws.Range("H99").Value = -1

$ws.Range("G100").Value = "Fallo"
$ws.Range("H100").Value = -1

$ws.Range("G101").Value = "Fallo"
$ws.Range("H101").Value = -1

$ws.Range("G102").Value = "Fallo"
$ws.Range("H102").Value = -1

# Append a new match row (row 105) at the end of the tracker
$ws.Range("A105").Value = 14601569

# Force the date-looking text to stay as plain text (avoid Excel auto
# converting "2025-09-13" into a date serial number/format)
$ws.Range("B105").NumberFormat = "@"
$ws.Range("B105").Value = "2025-09-13"
$ws.Range("B105").ClearFormats()

$ws.Range("C105").Value = "Nikola Bartunkova"
$ws.Range("D105").Value = "Magdalena Frech"
$ws.Range("E105").Value = "Gana Magdalena Frech"
$ws.Range("F105").Value = 1.5
